$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the notes text for the pareto_coefficient_2001_wid entry (cell K6)
$ws.Range("K6").Value = "Calculated from data by the World Inequality Databse. See downloadParetoCoefficient.R in the taxReform2001 folder for details."

# 2. Apply word-wrap style to K6 (same style as other "notes" cells, e.g. M1/M4) and set row height to 45
$ws.Range("K6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 45

# 3. Update the active selection in the sheet view to K6 instead of L6
$ws.Range("K6").Select()
